$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph (footer line that is being
# removed along with the blank paragraph before it and the copyright
# paragraph after it).
$findRange = $d.Content.Duplicate
$found = $findRange.Find.Execute(
    "Ver no Jupiter Salvar em pdf Salvar em docx",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $idx = $findRange.Paragraphs.Item(1).Index

    $prevPara = $d.Paragraphs.Item($idx - 1)
    $nextPara = $d.Paragraphs.Item($idx + 1)

    # Remove the blank paragraph, the "Ver no Jupiter..." paragraph and the
    # "© 2020 ..." paragraph that follows it, all in one shot (this also
    # removes the paragraph marks, so the surrounding paragraphs close up
    # correctly).
    $delRange = $d.Range($prevPara.Range.Start, $nextPara.Range.End)
    $delRange.Delete()
}
